$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.035827
$ws.Range("H2").Value = 3.107481
$ws.Range("I2").Value = 0.7561862865052227
$ws.Range("J2").Value = 0.7561862865052227
$ws.Range("M2").Value = 0.5134303333333333
$ws.Range("N2").Value = 1.540291
$ws.Range("O2").Value = 0.03326489761800302
$ws.Range("P2").Value = 0.03326489761800301
$ws.Range("Q2").Value = 0.5318250018856667
$ws.Range("R2").Value = 4.786425016970999
$ws.Range("S2").Value = 0.02515445940073413
$ws.Range("T2").Value = 0.02515445940073412
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.035827
$ws.Range("H3").Value = 3.107481
$ws.Range("I3").Value = 0.7561862865052227
$ws.Range("J3").Value = 0.7561862865052227
$ws.Range("M3").Value = 1.626140333333333
$ws.Range("N3").Value = 4.878420999999999
$ws.Range("O3").Value = 0.1053568287437347
$ws.Range("P3").Value = 0.1053568287437347
$ws.Range("Q3").Value = 1.684400063055667
$ws.Range("R3").Value = 15.159600567501
$ws.Range("S3").Value = 0.07966938908569146
$ws.Range("T3").Value = 0.07966938908569146
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.035827
$ws.Range("H4").Value = 3.107481
$ws.Range("I4").Value = 0.7561862865052227
$ws.Range("J4").Value = 0.7561862865052227
$ws.Range("M4").Value = 10.254745
$ws.Range("N4").Value = 30.764235
$ws.Range("O4").Value = 0.6643998618255804
$ws.Range("P4").Value = 0.6643998618255803
$ws.Range("Q4").Value = 10.622141749115
$ws.Range("R4").Value = 95.59927574203499
$ws.Range("S4").Value = 0.5024100642684687
$ws.Range("T4").Value = 0.5024100642684687
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.035827
$ws.Range("H5").Value = 3.107481
$ws.Range("I5").Value = 0.7561862865052227
$ws.Range("J5").Value = 0.7561862865052227
$ws.Range("M5").Value = 3.040282666666667
$ws.Range("N5").Value = 9.120848000000001
$ws.Range("O5").Value = 0.1969784118126819
$ws.Range("P5").Value = 0.1969784118126819
$ws.Range("Q5").Value = 3.149206873765333
$ws.Range("R5").Value = 28.342861863888
$ws.Range("S5").Value = 0.1489523737503284
$ws.Range("T5").Value = 0.1489523737503284
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.333977
$ws.Range("H6").Value = 1.001931
$ws.Range("I6").Value = 0.2438137134947773
$ws.Range("J6").Value = 0.2438137134947773
$ws.Range("M6").Value = 0.5134303333333333
$ws.Range("N6").Value = 1.540291
$ws.Range("O6").Value = 0.03326489761800302
$ws.Range("P6").Value = 0.03326489761800301
$ws.Range("Q6").Value = 0.1714739224356666
$ws.Range("R6").Value = 1.543265301921
$ws.Range("S6").Value = 0.008110438217268888
$ws.Range("T6").Value = 0.008110438217268888
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.333977
$ws.Range("H7").Value = 1.001931
$ws.Range("I7").Value = 0.2438137134947773
$ws.Range("J7").Value = 0.2438137134947773
$ws.Range("M7").Value = 1.626140333333333
$ws.Range("N7").Value = 4.878420999999999
$ws.Range("O7").Value = 0.1053568287437347
$ws.Range("P7").Value = 0.1053568287437347
$ws.Range("Q7").Value = 0.5430934701056666
$ws.Range("R7").Value = 4.887841230950999
$ws.Range("S7").Value = 0.02568743965804325
$ws.Range("T7").Value = 0.02568743965804326
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.333977
$ws.Range("H8").Value = 1.001931
$ws.Range("I8").Value = 0.2438137134947773
$ws.Range("J8").Value = 0.2438137134947773
$ws.Range("M8").Value = 10.254745
$ws.Range("N8").Value = 30.764235
$ws.Range("O8").Value = 0.6643998618255804
$ws.Range("P8").Value = 0.6643998618255803
$ws.Range("Q8").Value = 3.424848970864999
$ws.Range("R8").Value = 30.823640737785
$ws.Range("S8").Value = 0.1619897975571117
$ws.Range("T8").Value = 0.1619897975571117
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.333977
$ws.Range("H9").Value = 1.001931
$ws.Range("I9").Value = 0.2438137134947773
$ws.Range("J9").Value = 0.2438137134947773
$ws.Range("M9").Value = 3.040282666666667
$ws.Range("N9").Value = 9.120848000000001
$ws.Range("O9").Value = 0.1969784118126819
$ws.Range("P9").Value = 0.1969784118126819
$ws.Range("Q9").Value = 1.015384484165333
$ws.Range("R9").Value = 9.138460357488
$ws.Range("S9").Value = 0.04802603806235348
$ws.Range("T9").Value = 0.04802603806235348
